$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '41.857.32'
$ws.Range('E2').Value = '  +2.56%  '

# Row 3
$ws.Range('D3').Value = '2.232.05'
$ws.Range('E3').Value = '  +0.80%  '

# Row 4
$ws.Range('E4').Value = '  -0.04%  '

# Row 5
$ws.Range('D5').Value = '''232.42'
$ws.Range('E5').Value = '  +1.19%  '

# Row 6
$ws.Range('D6').Value = '''0.624'
$ws.Range('E6').Value = '  -1.32%  '

# Row 7
$ws.Range('D7').Value = '''60.77'
$ws.Range('E7').Value = '  -5.81%  '

# Row 8
$ws.Range('E8').Value = '  -0.02%  '

# Row 9
$ws.Range('E9').Value = '  +0.72%  '

# Row 10
$ws.Range('D10').Value = '''58.08'
$ws.Range('E10').Value = '  -1.77%  '

# Row 11
$ws.Range('D11').Value = '''0.0906'
$ws.Range('E11').Value = '  +4.51%  '

# Row 12
$ws.Range('E12').Value = '  -0.12%  '

# Row 13
$ws.Range('D13').Value = '2.564.58'
$ws.Range('E13').Value = '  +0.81%  '

# Row 14
$ws.Range('D14').Value = '''15.73'
$ws.Range('E14').Value = '  -1.17%  '

# Row 15
$ws.Range('D15').Value = '''22.66'
$ws.Range('E15').Value = '  +1.88%  '

# Row 16
$ws.Range('D16').Value = '''0.804'
$ws.Range('E16').Value = '  -2.09%  '

# Row 17
$ws.Range('D17').Value = '''5.62'
$ws.Range('E17').Value = '  +0.17%  '

# Row 18
$ws.Range('D18').Value = '2.249.56'
$ws.Range('E18').Value = '  +1.85%  '

# Row 19
$ws.Range('D19').Value = '41.808.40'
$ws.Range('E19').Value = '  +2.76%  '

# Row 20
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').Value = '0.0₃0907'
$ws.Range('E20').Value = '  +0.49%  '

# Row 21
$ws.Range('B21').Value = 'Litecoin'
$ws.Range('C21').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D21').Value = '''72.54'
$ws.Range('E21').Value = '  -2.00%  '

# Row 22
$ws.Range('D22').Value = '''6.13'
$ws.Range('E22').Value = '  -0.70%  '

# Row 23
$ws.Range('D23').Value = '''247.96'
$ws.Range('E23').Value = '  -1.02%  '

# Row 24
$ws.Range('E24').Value = '  -0.09%  '

# Row 25
$ws.Range('E25').Value = '  +0.11%  '

# Row 26
$ws.Range('E26').Value = '  +0.23%  '

# Row 27
$ws.Range('E27').Value = '  -0.23%  '

# Row 28
$ws.Range('B28').Value = 'Monero'
$ws.Range('C28').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D28').Value = '''169.63'
$ws.Range('E28').Value = '  -1.98%  '

# Row 29
$ws.Range('B29').Value = 'Kaspa'
$ws.Range('C29').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D29').Value = '''0.143'
$ws.Range('E29').Value = '  -0.56%  '

# Row 30
$ws.Range('D30').Value = '''19.95'
$ws.Range('E30').Value = '  -1.90%  '

# Row 31
$ws.Range('D31').Value = '''1.40'
$ws.Range('E31').Value = '  -2.59%  '

# Row 32
$ws.Range('D32').Value = '''2.66'
$ws.Range('E32').Value = '  -5.69%  '

# Row 33
$ws.Range('E33').Value = '  -1.35%  '

# Row 34
$ws.Range('D34').Value = '''5.06'
$ws.Range('E34').Value = '  +5.95%  '

# Row 35
$ws.Range('D35').Value = '''4.70'
$ws.Range('E35').Value = '  +0.74%  '

# Row 36
$ws.Range('D36').Value = '''0.0654'
$ws.Range('E36').Value = '  +3.58%  '

# Row 38
$ws.Range('B38').Value = 'LidoDAOToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D38').Value = '''2.40'
$ws.Range('E38').Value = '  -2.62%  '

# Row 39
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').Value = '''3.63'
$ws.Range('E39').Value = '  -5.00%  '

# Row 40
$ws.Range('B40').Value = 'BinanceUSD'
$ws.Range('C40').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D40').Value = '''1.00'
$ws.Range('E40').Value = '  +0.30%  '

# Row 41
$ws.Range('B41').Value = 'TerraClassic'
$ws.Range('C41').Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$ws.Range('D41').Value = '''0.000237'
$ws.Range('E41').Value = '  +14.71%  '

# Row 42
$ws.Range('E42').Value = '  +3.51%  '

# Row 43
$ws.Range('E43').Value = '  +0.37%  '

# Row 44
$ws.Range('D44').Value = '''1.22'
$ws.Range('E44').Value = '  -0.78%  '

# Row 45
$ws.Range('D45').Value = '''98.73'
$ws.Range('E45').Value = '  -2.64%  '

# Row 46
$ws.Range('E46').Value = '  +1.80%  '

# Row 47
$ws.Range('D47').Value = '1.471.33'
$ws.Range('E47').Value = '  -2.55%  '

# Row 48
$ws.Range('E48').Value = '  -12.18%  '

# Row 49
$ws.Range('D49').Value = '''16.65'
$ws.Range('E49').Value = '  -4.00%  '

# Row 50
$ws.Range('B50').Value = 'HuobiToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D50').Value = '''2.77'
$ws.Range('E50').Value = '  -2.38%  '

# Row 51
$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D51').Value = '''2.27'
$ws.Range('E51').Value = '  +6.79%  '
